# Trade #106 closed at 2026-02-17 09:18:41 - unknown UNKNOWN +0.000%
#
# This script updates the "live_trading_results" workbook to record the
# closing of trade #106:
#   - Summary sheet: Total Trades and Win Rate % are refreshed.
#   - Strategy Status sheet: the MarketMaking row's Trades and Win Rate % are refreshed.
#   - All Trades / MarketMaking sheets: a new trade-log row (#107 / trade 106) is appended.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet totals
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 106      # Total Trades
$wsSummary.Range("B9").Value = 42.45    # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet (MarketMaking row = row 4)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 106       # Trades
$wsStatus.Range("G4").Value = 42.45     # Win Rate %

# ---------------------------------------------------------------------------
# 3) Append the new closed-trade row to both the "All Trades" and
#    "MarketMaking" trade logs (they mirror one another since MarketMaking
#    is currently the only strategy generating trades).
# ---------------------------------------------------------------------------
function Add-TradeLogRow107([object]$ws) {
    $row = 107

    $ws.Cells.Item($row, 1).Value = 106          # Trade #

    # The Date column holds a literal "YYYY-MM-DD" text value in this sheet
    # (not an actual date type), so force text formatting before assigning
    # it to stop the string from being auto-parsed into a date serial
    # number, then drop the temporary number format again.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17" # Date
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = "09:18:35"                            # Time
    $ws.Cells.Item($row, 4).Value = "MarketMaking"                        # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"                                # Side
    $ws.Cells.Item($row, 6).Value = 0.83                                  # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.83                                  # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"                              # Status
    $ws.Cells.Item($row, 9).Value = 0                                     # P&L %
    $ws.Cells.Item($row, 10).Value = 0                                    # P&L $
    $ws.Cells.Item($row, 11).Value = 100.12                               # Capital After
    $ws.Cells.Item($row, 12).Value = 0                                    # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                                    # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                                  # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"   # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"                        # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.13                                 # Duration (min)
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeLogRow107 $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeLogRow107 $wsMarketMaking
